$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy formatting from F:G (old D:E, now shifted) into new D:E so new cells match style
$ws.Range("F7:G102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New quarterly data for columns D (2018-12-31) and E (2018-09-30)
$newColData = @{
    7 = @(43465, 43373)
    8 = @(202700, 185100)
    9 = @("NA", "NA")
    10 = @("NA", "NA")
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(-1300, -1400)
    17 = @(88900, 40400)
    18 = @(113800, 144700)
    20 = @(-89300, -79500)
    21 = @(37800, 78500)
    22 = @(0, 0)
    23 = @(24500, 65200)
    24 = @(1000, 10400)
    25 = @(0, 0)
    26 = @(23500, 54800)
    27 = @(23500, 54800)
    28 = @(0, 0)
    29 = @(2000, 3000)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(89300, 79500)
    33 = @(25500, 57800)
    34 = @(0, 0)
    35 = @(25500, 57800)
    38 = @(43465, 43373)
    41 = @(645100, 348700)
    42 = @(1735800, 956600)
    43 = @(0, 0)
    44 = @(0, 0)
    45 = @(0, 0)
    46 = @(0, 0)
    47 = @(0, 0)
    48 = @(283900, 277100)
    49 = @(195900, 197200)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @("NA", "NA")
    53 = @(0, 0)
    54 = @(23351100, 21462400)
    57 = @(0, 0)
    58 = @(0, 0)
    59 = @(177700, 171500)
    60 = @(0, 0)
    61 = @(82700, 78500)
    62 = @(0, 0)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(21122600, 19258900)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(1488400, 1477700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(2228500, 2203500)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(25500, 57800)
    83 = @(13300, 13200)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(73000, 125300)
    91 = @(-21200, -14600)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-1210300, -457200)
    96 = @(-14800, -14500)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(1826700, 844200)
    101 = @(0, 0)
    102 = @(689400, 512300)
}

foreach ($r in $newColData.Keys) {
    $vals = $newColData[$r]
    $ws.Range("D$r").Value2 = $vals[0]
    $ws.Range("E$r").Value2 = $vals[1]
}

Write-Host "Done applying column insert and new data"